# Dashboard updated 6 cases and login page updated - Loginpage and dashboard page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Remove the old hyperlink (was on B2) before we start moving data around ----
foreach ($hl in $ws.Hyperlinks) {
  $hl.Delete()
}

# ---- Drop the old "Result" column (column D) entirely ----
$ws.Range("D1:D3").Clear() | Out-Null

# ---- Re-purpose the remaining 3 columns into: email / password / loginstatus ----
# Header row
$ws.Range("A1").Value2 = "email"
$ws.Range("B1").Value2 = "password"
$ws.Range("C1").Value2 = "loginstatus"

# Data rows
$ws.Range("A2").Value2 = "yatutor200@gmail.com"
$ws.Range("B2").Value2 = "India@2020"
$ws.Range("C2").Value2 = "Passed"

$ws.Range("A3").Value2 = "yatutor200@gmail.com"
$ws.Range("B3").Value2 = "India@2021"
$ws.Range("C3").Value2 = "Failed"

# The old hyperlink cell (B2) goes back to plain/default formatting
$ws.Range("B2").Style = "Normal"

# ---- Re-create the hyperlink, now anchored on A2 (the email cell) ----
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:yatutor200@gmail.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

# ---- Column widths for the new layout ----
$ws.Columns.Item(1).ColumnWidth = 19.666666666666664
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Columns.Item(3).ColumnWidth = 9.333333333333332

# ---- Selection moves to G11 ----
$ws.Range("G11").Select() | Out-Null
